# Update capital structure database for Estonia beverage (alcoholic) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

# Row 3: fix company name spelling ("AS" -> "aktsiaselts")
$ws.Range("B3").Value = "aktsiaselts Linda Nektar (TLSE:LINDA)"

# New columns D (historical_growth_revenue_last_5_years) and
# E (historical_growth_net_income_last_5_years) for both data rows.
$ws.Range("D2:D3").Value = -0.0196
$ws.Range("E2:E3").Value = -0.358

# Recomputed financial ratios for row 2 and row 3 (identical values).
foreach ($row in 2, 3) {
    $ws.Cells.Item($row, 7).Value  = 0.1484962406015038    # G - ebitdard_margin
    $ws.Cells.Item($row, 8).Value  = 0.1379699248120301    # H - ebitda_margin
    $ws.Cells.Item($row, 9).Value  = 0.01428571428571429   # I - operating_margin
    $ws.Cells.Item($row, 10).Value = 0.007962529274004683  # J - after_tax_operating_margin
    $ws.Cells.Item($row, 11).Value = 0.034                 # K - trailing_net_income
    $ws.Cells.Item($row, 12).Value = 0.01278195488721805   # L - net_margin
    $ws.Cells.Item($row, 13).Value = 0.142                 # M - cash_returned
    $ws.Cells.Item($row, 14).Value = 0.009793103448275862  # N - cash_returned_market_cap
    $ws.Cells.Item($row, 15).Value = 4.176470588235293     # O - cash_returned_net_income
    $ws.Cells.Item($row, 16).Value = 0.142                 # P - dividends
    $ws.Cells.Item($row, 17).Value = 0.009793103448275862  # Q - dividend_yield
    $ws.Cells.Item($row, 18).Value = 4.176470588235293     # R - payout_ratio

    $ws.Cells.Item($row, 21).Value = 0.784                 # U - cash
    $ws.Cells.Item($row, 22).Value = 0.05406896551724138   # V - cash_market_cap
    $ws.Cells.Item($row, 23).Value = 0.007962529274004685  # W - roe
    $ws.Cells.Item($row, 24).Value = 0.04895538255292466   # X - cost_equity
    $ws.Cells.Item($row, 25).Value = -0.04099285327891997  # Y - roe_cost_equity
    $ws.Cells.Item($row, 26).Value = 0.808264965056214     # Z - sales_invested_capital
    $ws.Cells.Item($row, 27).Value = 0.006435833445412477  # AA - roic
    $ws.Cells.Item($row, 28).Value = 0.04895538255292466   # AB - cost_capital
    $ws.Cells.Item($row, 29).Value = -0.04251954910751218  # AC - roic_cost_capital

    $ws.Cells.Item($row, 33).Value = -0.784                # AG - net_debt

    $ws.Cells.Item($row, 36).Value = -0.05715952172645086  # AJ - net_debt_market_capital
    $ws.Cells.Item($row, 37).Value = -0.2350119904076738   # AK - net_debt_book_capital

    $ws.Cells.Item($row, 42).Value = -1.708061002178649    # AP - net_debt_ebitda
}
